# Automatic data update 2026-02-27 23:50 -- refresh DATA_EXTRACCIO timestamps
# and the handful of re-measured values pulled in during that meteo.cat sync.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-27 23:48:24"
$ws.Range("O2").Value = "5.2 °C"
$ws.Range("E3").Value = "2026-02-27 23:48:27"
$ws.Range("E4").Value = "2026-02-27 23:48:29"
$ws.Range("E5").Value = "2026-02-27 23:48:32"
$ws.Range("H5").Value = "'49%"
$ws.Range("O5").Value = "4.3 °C"
$ws.Range("E6").Value = "2026-02-27 23:48:34"
$ws.Range("E7").Value = "2026-02-27 23:48:37"
$ws.Range("H7").Value = "'85%"
$ws.Range("E8").Value = "2026-02-27 23:48:39"
$ws.Range("E9").Value = "2026-02-27 23:48:42"
$ws.Range("E10").Value = "2026-02-27 23:48:44"
$ws.Range("O10").Value = "10.6 °C"
$ws.Range("E11").Value = "2026-02-27 23:48:45"
$ws.Range("O11").Value = "8.2 °C"
$ws.Range("E12").Value = "2026-02-27 23:48:46"
$ws.Range("O12").Value = "10.6 °C"
$ws.Range("E13").Value = "2026-02-27 23:48:47"
$ws.Range("E14").Value = "2026-02-27 23:48:48"
$ws.Range("H14").Value = "'93%"
$ws.Range("O14").Value = "10.1 °C"
$ws.Range("E15").Value = "2026-02-27 23:48:50"
$ws.Range("H15").Value = "'89%"
$ws.Range("N15").Value = "6.4 °C 23:18 TU"
$ws.Range("O15").Value = "10.6 °C"
$ws.Range("E16").Value = "2026-02-27 23:48:51"
$ws.Range("H16").Value = "'44%"
$ws.Range("N16").Value = "-0.8 °C 23:02 TU"
$ws.Range("O16").Value = "2.4 °C"
$ws.Range("E17").Value = "2026-02-27 23:48:52"
$ws.Range("E18").Value = "2026-02-27 23:48:53"
$ws.Range("N18").Value = "7.9 °C 23:23 TU"
$ws.Range("E19").Value = "2026-02-27 23:48:54"
$ws.Range("N19").Value = "6.2 °C 23:12 TU"
$ws.Range("O19").Value = "10.1 °C"
$ws.Range("E20").Value = "2026-02-27 23:48:55"
$ws.Range("K20").Value = "16.7 MJ/m2"
$ws.Range("O20").Value = "3.0 °C"
$ws.Range("E21").Value = "2026-02-27 23:48:56"
$ws.Range("H21").Value = "'60%"
$ws.Range("E22").Value = "2026-02-27 23:48:59"
$ws.Range("H22").Value = "'50%"
$ws.Range("N22").Value = "-0.9 °C 23:15 TU"
$ws.Range("E23").Value = "2026-02-27 23:49:01"
$ws.Range("N23").Value = "0.7 °C 23:21 TU"
$ws.Range("O23").Value = "3.4 °C"
$ws.Range("E24").Value = "2026-02-27 23:49:04"
$ws.Range("L24").Value = "22.7 km/h - 161º 23:22 TU"
$ws.Range("E25").Value = "2026-02-27 23:49:06"
$ws.Range("H25").Value = "'36%"
$ws.Range("E26").Value = "2026-02-27 23:49:09"
$ws.Range("K26").Value = "16.0 MJ/m2"
$ws.Range("E27").Value = "2026-02-27 23:49:11"
$ws.Range("E28").Value = "2026-02-27 23:49:14"
$ws.Range("E29").Value = "2026-02-27 23:49:16"
$ws.Range("H29").Value = "'88%"
$ws.Range("E30").Value = "2026-02-27 23:49:18"
$ws.Range("N30").Value = "8.2 °C 23:28 TU"
$ws.Range("E31").Value = "2026-02-27 23:49:20"
$ws.Range("O31").Value = "10.3 °C"
$ws.Range("E32").Value = "2026-02-27 23:49:23"
$ws.Range("H32").Value = "'53%"
$ws.Range("E33").Value = "2026-02-27 23:49:26"
$ws.Range("J33").Value = "1023.5 hPa"
$ws.Range("E34").Value = "2026-02-27 23:49:28"
$ws.Range("O34").Value = "4.3 °C"
$ws.Range("E35").Value = "2026-02-27 23:49:31"
$ws.Range("H35").Value = "'43%"
$ws.Range("O35").Value = "11.7 °C"
$ws.Range("E36").Value = "2026-02-27 23:49:33"
$ws.Range("H36").Value = "'92%"
$ws.Range("E37").Value = "2026-02-27 23:49:36"
$ws.Range("E38").Value = "2026-02-27 23:49:38"
$ws.Range("E39").Value = "2026-02-27 23:49:41"
$ws.Range("H39").Value = "'33%"
$ws.Range("O39").Value = "4.3 °C"
$ws.Range("E40").Value = "2026-02-27 23:49:43"
$ws.Range("H40").Value = "'69%"
$ws.Range("O40").Value = "8.7 °C"
$ws.Range("E41").Value = "2026-02-27 23:49:45"
$ws.Range("H41").Value = "'84%"
$ws.Range("J41").Value = "1024.4 hPa"
$ws.Range("O41").Value = "11.2 °C"
$ws.Range("E42").Value = "2026-02-27 23:49:48"
$ws.Range("O42").Value = "11.1 °C"
$ws.Range("E43").Value = "2026-02-27 23:49:50"
$ws.Range("O43").Value = "9.1 °C"
$ws.Range("E44").Value = "2026-02-27 23:49:53"
$ws.Range("K44").Value = "17.0 MJ/m2"
$ws.Range("O44").Value = "1.8 °C"
$ws.Range("E45").Value = "2026-02-27 23:49:55"
$ws.Range("H45").Value = "'48%"
$ws.Range("J45").Value = "1022.0 hPa"
$ws.Range("O45").Value = "11.5 °C"
$ws.Range("E46").Value = "2026-02-27 23:49:57"
$ws.Range("K46").Value = "13.1 MJ/m2"
